$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New data rows (276-293) -------------------------------------------------
# Columns: A season, B Session(date), C game, D loc, F player, G score,
# H special, I place, J position, K desert, L sum(formula), M count_num,
# N count_res, O..V p_sum_*, W..AB port_*, AC..AH prox_port_*

$rows = @(
    @(2024,45646,40,"Husum (Memeler Str.)","JHC",8,2,3,"first","outer",5,4,5,0,3,5,8,2,4,4,0,0,0,0,0,0,0,0,0,0,1,0),
    @(2024,45646,40,"Husum (Memeler Str.)","PF",10,2,2,"second","outer",6,4,5,17,4,0,1,0,1,4,0,0,0,0,0,0,0,1,0,0,0,0),
    @(2024,45646,40,"Husum (Memeler Str.)","MF",13,5,1,"third","outer",6,5,2,3,5,2,8,0,8,5,0,0,0,0,0,0,0,0,0,0,1,0),
    @(2024,45646,41,"Husum (Memeler Str.)","PF",12,3,2,"first","outer",6,5,5,6,5,9,4,5,0,0,0,0,0,0,0,0,0,0,0,0,0,1),
    @(2024,45646,41,"Husum (Memeler Str.)","JHC",9,2,3,"second","outer",6,5,4,6,3,5,9,0,9,0,0,0,0,0,0,0,0,0,0,0,0,0),
    @(2024,45646,41,"Husum (Memeler Str.)","MF",13,5,1,"third","outer",5,5,3,3,5,7,4,3,4,5,0,0,0,0,0,0,0,1,1,0,0,0),
    @(2024,45646,42,"Husum (Memeler Str.)","JHC",13,6,1,"first","outer",5,4,3,22,0,4,5,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0),
    @(2024,45646,42,"Husum (Memeler Str.)","PF",7,2,2,"second","outer",6,4,5,0,5,3,9,1,9,0,0,0,0,0,0,0,0,0,1,0,1,1),
    @(2024,45646,42,"Husum (Memeler Str.)","MF",6,2,3,"third","outer",5,5,4,3,6,4,5,4,5,2,0,0,0,0,0,0,0,0,0,0,0,0),
    @(2024,45646,43,"Husum (Memeler Str.)","JHC",8,2,2,"first","outer",5,3,0,27,5,4,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0),
    @(2024,45646,43,"Husum (Memeler Str.)","PF",13,8,1,"second","outer",6,4,5,0,5,3,8,5,5,0,0,0,0,0,0,0,1,0,1,0,0,0),
    @(2024,45646,43,"Husum (Memeler Str.)","MF",8,1,3,"third","outer",5,5,7,8,4,4,3,0,3,4,0,0,0,0,0,0,0,0,0,0,0,0),
    @(2024,45646,44,"Husum (Memeler Str.)","JHC",8,2,2,"first","outer",5,4,4,13,4,10,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,1),
    @(2024,45646,44,"Husum (Memeler Str.)","MF",5,1,3,"second","outer",6,5,2,6,5,3,8,2,0,5,0,0,0,0,0,0,0,0,0,0,0,0),
    @(2024,45646,44,"Husum (Memeler Str.)","PF",13,5,1,"third","outer",4,4,9,3,0,2,9,9,0,0,0,0,0,0,0,0,0,1,0,0,0,0),
    @(2024,45646,45,"Husum (Memeler Str.)","JHC",13,5,1,"first","outer",5,5,4,8,3,5,4,0,1,3,0,0,0,0,0,0,0,0,0,0,0,0),
    @(2024,45646,45,"Husum (Memeler Str.)","MF",10,2,2,"second","outer",5,5,2,5,3,12,3,2,0,3,0,0,0,0,0,0,0,0,0,0,0,0),
    @(2024,45646,45,"Husum (Memeler Str.)","PF",9,2,3,"third","outer",5,5,5,3,6,5,2,5,2,3,0,0,0,0,0,0,0,1,0,0,0,0)
)

$startRow = 276
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]

    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("F$r").Value = $row[4]
    $ws.Range("G$r").Value = $row[5]
    $ws.Range("H$r").Value = $row[6]
    $ws.Range("I$r").Value = $row[7]
    $ws.Range("J$r").Value = $row[8]
    $ws.Range("K$r").Value = $row[9]
    $ws.Range("L$r").Formula = "=SUM(O$r`:V$r)"
    $ws.Range("M$r").Value = $row[10]
    $ws.Range("N$r").Value = $row[11]
    $ws.Range("O$r").Value = $row[12]
    $ws.Range("P$r").Value = $row[13]
    $ws.Range("Q$r").Value = $row[14]
    $ws.Range("R$r").Value = $row[15]
    $ws.Range("S$r").Value = $row[16]
    $ws.Range("T$r").Value = $row[17]
    $ws.Range("U$r").Value = $row[18]
    $ws.Range("V$r").Value = $row[19]
    $ws.Range("W$r").Value = $row[20]
    $ws.Range("X$r").Value = $row[21]
    $ws.Range("Y$r").Value = $row[22]
    $ws.Range("Z$r").Value = $row[23]
    $ws.Range("AA$r").Value = $row[24]
    $ws.Range("AB$r").Value = $row[25]
    $ws.Range("AC$r").Value = $row[26]
    $ws.Range("AD$r").Value = $row[27]
    $ws.Range("AE$r").Value = $row[28]
    $ws.Range("AF$r").Value = $row[29]
    $ws.Range("AG$r").Value = $row[30]
    $ws.Range("AH$r").Value = $row[31]
}

# Apply the same date style (column B) used by existing rows onto the new ones
$ws.Range("B275").Copy()
$ws.Range("B276:B293").PasteSpecial(-4122)

# --- Corrections to three existing score values (G column) -----------------
$ws.Range("G27").Value = 13
$ws.Range("G172").Value = 13
$ws.Range("G258").Value = 13

# --- View / navigation state updates ---------------------------------------
$ws.Range("B275").Select()
$activeWin = $ws.Application.ActiveWindow
$activeWin.ScrollRow = 272
$ws.Range("E288").Select()

# --- AutoFilter / defined name range extension ------------------------------
$ws.Range("A2:AP290").AutoFilter()
